$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-12 -> 2023-09-13) for every data row (rows 2-219).
for ($r = 2; $r -le 219; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}
